# Generate Report for Archive
#
# 1. Update the localization status text from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview, zh-cn, de-de sheets).
# 2. Narrow the now-shorter "Status"/"zh-cn"/"de-de" columns to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text (was "Ready for handoff") -------------------------
$overview.Range("E2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# --- Shrink the affected columns to match the new, shorter text -----------
# (target stored width ~= 13.41 "Excel width" units; ColumnWidth is specified
#  in un-padded character units, so 12.5 here serializes to the intended
#  stored column width in the saved OOXML)
$overview.Columns("E").ColumnWidth = 12.5
$overview.Columns("F").ColumnWidth = 12.5
$zhcn.Columns("C").ColumnWidth = 12.5
$dede.Columns("C").ColumnWidth = 12.5
